$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string values (B1 = "1:4" -> "1:4 Conditioned", C1 = "4:1" -> "4:1 Conditioned")
$ws.Range("B1").Value = "1:4 Conditioned"
$ws.Range("C1").Value = "4:1 Conditioned"

# Update selection to E10
$ws.Range("E10").Select()

# Set column widths for columns B and C (~19.83 and ~16.16 characters wide)
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 15.3
